# Commit: "cap nhat anh cho 20 san pham" (update images for 20 products)
#
# The 20 products that lost their extra product-gallery images (pic2 / pic3
# columns, i.e. worksheet columns E / F) now reference the placeholder
# string "null" instead of the old "CuaHangTrangSuc\productsInfo\products2\*.png"
# / "...\products3\*.png" paths. Writing the literal string content (rather
# than poking shared-string indices) lets Excel's own shared-strings
# bookkeeping drop the now-unreferenced picture-path strings and append the
# single new "null" string, exactly like a real save would.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F2").Value = "null"

$ws.Range("E4").Value = "null"
$ws.Range("F4").Value = "null"

$ws.Range("E7").Value = "null"

$ws.Range("E10").Value = "null"
$ws.Range("F10").Value = "null"

$ws.Range("E11").Value = "null"
$ws.Range("F11").Value = "null"

$ws.Range("F12").Value = "null"

$ws.Range("E13").Value = "null"
$ws.Range("F13").Value = "null"

$ws.Range("E14").Value = "null"
$ws.Range("F14").Value = "null"

$ws.Range("F15").Value = "null"

$ws.Range("E16").Value = "null"
$ws.Range("F16").Value = "null"

$ws.Range("E17").Value = "null"
$ws.Range("F17").Value = "null"

$ws.Range("E18").Value = "null"
$ws.Range("F18").Value = "null"

$ws.Range("E19").Value = "null"
$ws.Range("F19").Value = "null"

$ws.Range("E20").Value = "null"
$ws.Range("F20").Value = "null"

$ws.Range("E21").Value = "null"
$ws.Range("F21").Value = "null"

# Restore the selection left on the sheet when the author saved.
$ws.Activate()
$ws.Range("E22").Select()
